$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new consolidated sell data
$ws.Range("B2").Value = "1AYB-5AYB-1AP-1M-4AYB"
$ws.Range("C2").Value = "1-1-1-1-1"
$ws.Range("D2").Value = 32100
$ws.Range("E2").Value = 45815.80757366943

# Remove the now-obsolete rows 3 through 6
$ws.Range("A3:E6").Delete()
